# Refresh the cryptos price/volume snapshot (+ restore two swapped row pairs).
# D-column numeric-looking values are written with a leading quote-prefix (")
# so Excel keeps them as text (matching the source's dotted/padded formatting)
# instead of auto-coercing them into Number cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '56.703.90'
$ws.Range('E2').Value = '  +3.20%  '
# Row 3
$ws.Range('D3').Value = '2.325.77'
$ws.Range('E3').Value = '  +2.30%  '
# Row 4
$ws.Range('E4').Value = '  +0.04%  '
# Row 5
$ws.Range('D5').Value = "'517.42"
$ws.Range('E5').Value = '  +2.26%  '
# Row 6
$ws.Range('D6').Value = "'135.79"
$ws.Range('E6').Value = '  +5.85%  '
# Row 7
$ws.Range('D7').Value = "'0.995"
$ws.Range('E7').Value = '  -0.07%  '
# Row 8
$ws.Range('D8').Value = "'0.538"
$ws.Range('E8').Value = '  +1.69%  '
# Row 9
$ws.Range('D9').Value = '2.345.42'
$ws.Range('E9').Value = '  +2.76%  '
# Row 10
$ws.Range('D10').Value = "'0.103"
$ws.Range('E10').Value = '  +4.20%  '
# Row 11
$ws.Range('E11').Value = '  -1.05%  '
# Row 12
$ws.Range('D12').Value = "'5.35"
$ws.Range('E12').Value = '  +5.20%  '
# Row 13
$ws.Range('D13').Value = "'0.343"
$ws.Range('E13').Value = '  +0.34%  '
# Row 14
$ws.Range('D14').Value = "'24.04"
$ws.Range('E14').Value = '  +2.00%  '
# Row 15
$ws.Range('D15').Value = '2.742.88'
$ws.Range('E15').Value = '  +2.43%  '
# Row 16
$ws.Range('D16').Value = '56.737.76'
$ws.Range('E16').Value = '  +3.21%  '
# Row 17
$ws.Range('D17').Value = "'0.0000135"
$ws.Range('E17').Value = '  +2.96%  '
# Row 18
$ws.Range('D18').Value = '2.356.68'
$ws.Range('E18').Value = '  +4.08%  '
# Row 19
$ws.Range('D19').Value = "'10.56"
$ws.Range('E19').Value = '  +1.89%  '
# Row 20
$ws.Range('B20').Value = 'BitcoinCash'
$ws.Range('C20').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D20').Value = "'327.02"
$ws.Range('E20').Value = '  +3.92%  '
# Row 21
$ws.Range('B21').Value = 'Polkadot'
$ws.Range('C21').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D21').Value = "'4.23"
$ws.Range('E21').Value = '  +0.87%  '
# Row 22
$ws.Range('D22').Value = "'6.59"
$ws.Range('E22').Value = '  +0.27%  '
# Row 23
$ws.Range('E23').Value = '  +0.19%  '
# Row 24
$ws.Range('D24').Value = "'60.83"
$ws.Range('E24').Value = '  +1.68%  '
# Row 25
$ws.Range('E25').Value = '  +6.78%  '
# Row 26
$ws.Range('D26').Value = "'0.992"
$ws.Range('E26').Value = '  -0.36%  '
# Row 27
$ws.Range('D27').Value = "'8.01"
$ws.Range('E27').Value = '  +6.23%  '
# Row 28
$ws.Range('D28').Value = "'1.29"
$ws.Range('E28').Value = '  +11.70%  '
# Row 29
$ws.Range('D29').Value = '0.0₃0744'
$ws.Range('E29').Value = '  +5.52%  '
# Row 30
$ws.Range('D30').Value = "'168.22"
$ws.Range('E30').Value = '  -1.67%  '
# Row 31
$ws.Range('D31').Value = "'1.70"
$ws.Range('E31').Value = '  +3.60%  '
# Row 32
$ws.Range('D32').Value = "'6.22"
$ws.Range('E32').Value = '  +1.10%  '
# Row 33
$ws.Range('D33').Value = "'18.51"
$ws.Range('E33').Value = '  +2.98%  '
# Row 34
$ws.Range('E34').Value = '  -0.02%  '
# Row 35
$ws.Range('D35').Value = "'0.994"
$ws.Range('E35').Value = '  +0.02%  '
# Row 37
$ws.Range('D37').Value = "'0.922"
# Row 38
$ws.Range('D38').Value = "'4.02"
$ws.Range('E38').Value = '  +3.16%  '
# Row 39
$ws.Range('D39').Value = "'1.56"
$ws.Range('E39').Value = '  +7.16%  '
# Row 40
$ws.Range('D40').Value = "'38.36"
$ws.Range('E40').Value = '  +4.49%  '
# Row 41
$ws.Range('B41').Value = 'Aave'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D41').Value = "'142.55"
$ws.Range('E41').Value = '  +4.45%  '
# Row 42
$ws.Range('B42').Value = 'PolygonEcosystemToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D42').Value = "'0.380"
$ws.Range('E42').Value = '  +1.58%  '
# Row 43
$ws.Range('D43').Value = "'3.61"
$ws.Range('E43').Value = '  +3.87%  '
# Row 44
$ws.Range('D44').Value = "'5.26"
$ws.Range('E44').Value = '  +7.82%  '
# Row 45
$ws.Range('D45').Value = "'279.09"
$ws.Range('E45').Value = '  +7.95%  '
# Row 46
$ws.Range('D46').Value = "'0.0937"
$ws.Range('E46').Value = '  +1.86%  '
# Row 47
$ws.Range('D47').Value = "'0.0508"
$ws.Range('E47').Value = '  +0.44%  '
# Row 48
$ws.Range('D48').Value = "'0.563"
$ws.Range('E48').Value = '  +2.86%  '
# Row 49
$ws.Range('E49').Value = '  +2.94%  '
# Row 50
$ws.Range('D50').Value = "'17.88"
$ws.Range('E50').Value = '  +9.00%  '
# Row 51
$ws.Range('E51').Value = '  +1.75%  '
